$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.508.10'
$ws.Range("E2").Value = '  -2.20%  '
$ws.Range("D3").Value = '2.482.94'
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.36'
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.08'
$ws.Range("E6").Value = '  -4.66%  '
$ws.Range("E7").Value = '  -2.07%  '
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("E9").Value = '  -3.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.49'
$ws.Range("E10").Value = '  -4.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0781'
$ws.Range("E11").Value = '  -2.22%  '
$ws.Range("E12").Value = '  +0.25%  '
$ws.Range("D13").Value = '2.866.93'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  -3.78%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.557.46'
$ws.Range("E15").Value = '  +2.27%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.58'
$ws.Range("E16").Value = '  +1.42%  '
$ws.Range("E17").Value = '  -1.62%  '
$ws.Range("D18").Value = '41.456.10'
$ws.Range("E18").Value = '  -2.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.34'
$ws.Range("E19").Value = '  -3.55%  '
$ws.Range("D20").Value = '0.0₃0925'
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.21'
$ws.Range("E21").Value = '  -6.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.99'
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.09'
$ws.Range("E24").Value = '  -2.62%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  -3.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.08'
$ws.Range("E27").Value = '  -4.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.83'
$ws.Range("E29").Value = '  -1.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.22'
$ws.Range("E30").Value = '  -1.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.84'
$ws.Range("E31").Value = '  -2.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.49'
$ws.Range("E32").Value = '  -5.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.58'
$ws.Range("E33").Value = '  -2.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.25'
$ws.Range("E34").Value = '  +5.75%  '
$ws.Range("E35").Value = '  -3.27%  '
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.45'
$ws.Range("E37").Value = '  -12.17%  '
$ws.Range("E38").Value = '  -3.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.115'
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.101'
$ws.Range("E40").Value = '  -4.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.14'
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.79'
$ws.Range("E43").Value = '  -6.82%  '
$ws.Range("D44").Value = '1.985.95'
$ws.Range("E44").Value = '  -1.08%  '
$ws.Range("E45").Value = '  -2.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.03'
$ws.Range("E46").Value = '  -6.21%  '
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("D48").Value = '2.729.72'
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '69.51'
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '97.30'
$ws.Range("E50").Value = '  -2.56%  '
$ws.Range("E51").Value = '  -4.69%  '

$numericTextCells = @("D5","D6","D10","D11","D16","D19","D21","D22","D23","D27","D28","D29","D30","D31","D32","D33","D34","D37","D39","D40","D41","D43","D46","D49","D50")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).ClearFormats()
}
